$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 212, shifting rows 212:654 down to 213:655
$ws.Rows.Item(212).Insert()

$ws.Range("A212").NumberFormat = "@"
$ws.Range("A212").Value = "125"
$ws.Range("B212").Value = "Patrick Mckee , (Minnesota) F Killian Cardinale , (West Virginia), 3:25"
